# Auto commit at 2025-10-13  7:47:43.56
# Appends two new daily rows (2025-10-12) for the two charging stations
# to the bottom of the data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 84 - 四方坪站充电量(kw) for date serial 45942 (2025-10-12)
$ws.Range("A84").Value = 45942
$ws.Range("B84").Value = "四方坪站充电量(kw)"
$ws.Range("C84").Value = 909.57299999999998
$ws.Range("D84").Value = 961.38400000000001
$ws.Range("E84").Value = 413.96199999999993
$ws.Range("F84").Value = 575.78
$ws.Range("G84").Value = 323.08099999999996
$ws.Range("H84").Value = 620.75
$ws.Range("I84").Value = 395.74799999999993
$ws.Range("J84").Value = 198.625
$ws.Range("K84").Value = 209.88399999999999
$ws.Range("L84").Value = 219.06999999999996
$ws.Range("M84").Value = 292.64300000000003
$ws.Range("N84").Value = 58.23
$ws.Range("O84").Value = 832.62799999999982
$ws.Range("P84").Value = 1288.6610000000001
$ws.Range("Q84").Value = 505.65899999999993
$ws.Range("R84").Value = 351.233
$ws.Range("S84").Value = 396.98199999999997
$ws.Range("T84").Value = 232.09
$ws.Range("U84").Value = 54.519999999999996
$ws.Range("V84").Value = 85.179999999999993
$ws.Range("W84").Value = 82.22
$ws.Range("X84").Value = 46.8
$ws.Range("Y84").Value = 7.78
$ws.Range("Z84").Value = 162.553

# Row 85 - 高岭站充电量(kw) for date serial 45942 (2025-10-12)
$ws.Range("A85").Value = 45942
$ws.Range("B85").Value = "高岭站充电量(kw)"
$ws.Range("C85").Value = 426.35700000000008
$ws.Range("D85").Value = 395.31300000000005
$ws.Range("E85").Value = 216.524
$ws.Range("F85").Value = 34.25
$ws.Range("G85").Value = 40.296999999999997
$ws.Range("H85").Value = 185.99999999999997
$ws.Range("I85").Value = 101.64100000000002
$ws.Range("J85").Value = 212.07299999999998
$ws.Range("K85").Value = 294.39600000000002
$ws.Range("L85").Value = 87.300999999999988
$ws.Range("M85").Value = 112.85899999999999
$ws.Range("N85").Value = 120.09100000000001
$ws.Range("O85").Value = 396.44400000000002
$ws.Range("P85").Value = 368.49200000000002
$ws.Range("Q85").Value = 159.042
$ws.Range("R85").Value = 214.59900000000002
$ws.Range("S85").Value = 72.312999999999988
$ws.Range("T85").Value = 105.23
$ws.Range("U85").Value = 51.852999999999994
$ws.Range("V85").Value = 51.391999999999996
$ws.Range("W85").Value = 61.811999999999998
$ws.Range("X85").Value = 28.067
$ws.Range("Y85").Value = 0
$ws.Range("Z85").Value = 41.29

# Move/show the selection the way the author left it after entering the data
$ws.Range("F89").Select()
